# Rename 'Codelists' sheet to 'Cells' and update the active sheet/selection
# state to match (Close #256).

$wb = $excel.ActiveWorkbook

# Rename the "Codelists" sheet to "Cells".
$ws = $wb.Worksheets.Item("Codelists")
$ws.Name = "Cells"

# Restore the selection on the "Table" sheet (no longer the active tab).
$tableWs = $wb.Worksheets.Item("Table")
$tableWs.Range("L18").Select()

# Make the renamed "Cells" sheet the active sheet/tab, with its own selection.
$ws.Activate()
$ws.Range("I6").Select()
